# Generate Report for Handoff
#
# The two localization entries ("0d2067ae-2ca4-448f-bc8e-89192d7c768c" and
# "e3b61664-96dc-4ab8-bc89-9c0d7fefc835") swap row order on every sheet, and
# the 0d2067ae entry moves from "Handed back: in sync with en-US" to
# "Ready for handoff" with fresh handoff timestamps.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "2016-03-23 18:55:52"

$ws.Range("A3").Value = "0d2067ae-2ca4-448f-bc8e-89192d7c768c.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-03-23 18:57:33"

$ws.Hyperlinks.Delete()
$null = $ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/91dfa5fecd8eb92d70a2f5fb7a12f3cf15aa0586/e2e/0d2067ae-2ca4-448f-bc8e-89192d7c768c.md", "", "", "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.md")
$null = $ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/91dfa5fecd8eb92d70a2f5fb7a12f3cf15aa0586/e2e/e3b61664-96dc-4ab8-bc89-9c0d7fefc835.md", "", "", "0d2067ae-2ca4-448f-bc8e-89192d7c768c.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.5b909c6e9d6b4335b3af5e7cc35d338fdda5ef7d.zh-cn.xlf"
$ws.Range("E2").Value = "2016-03-23 18:55:47"
$ws.Range("F2").Value = "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.md"
$ws.Range("G2").Value = "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.5b909c6e9d6b4335b3af5e7cc35d338fdda5ef7d.zh-cn.xlf"
$ws.Range("H2").Value = "2016-03-23 18:56:33"
$ws.Range("J2").Value = "Include"

$ws.Range("A3").Value = "0d2067ae-2ca4-448f-bc8e-89192d7c768c.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "0d2067ae-2ca4-448f-bc8e-89192d7c768c.fdb0e26b708f7757c927665c32014118d81c9a82.zh-cn.xlf"
$ws.Range("E3").Value = "2016-03-23 18:57:29"
$ws.Range("F3").Value = "0d2067ae-2ca4-448f-bc8e-89192d7c768c.md"
$ws.Range("G3").Value = "0d2067ae-2ca4-448f-bc8e-89192d7c768c.fdb0e26b708f7757c927665c32014118d81c9a82.zh-cn.xlf"
$ws.Range("H3").Value = "2016-03-23 18:56:33"
$ws.Range("J3").Value = "Include"

$ws.Hyperlinks.Delete()
$null = $ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/91dfa5fecd8eb92d70a2f5fb7a12f3cf15aa0586/e2e/0d2067ae-2ca4-448f-bc8e-89192d7c768c.md", "", "", "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.md")
$null = $ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0f75b98accff473613980d6d00c49f8e386eb314/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0d2067ae-2ca4-448f-bc8e-89192d7c768c.fdb0e26b708f7757c927665c32014118d81c9a82.zh-cn.xlf", "", "", "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.5b909c6e9d6b4335b3af5e7cc35d338fdda5ef7d.zh-cn.xlf")
$null = $ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/30fdea63b4e1cf4e0d65702035e6c8e5d3bb2d6f/e2e/0d2067ae-2ca4-448f-bc8e-89192d7c768c.md", "", "", "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.md")
$null = $ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3436c8a9e40370fec3686ff289c357d4498874c7/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0d2067ae-2ca4-448f-bc8e-89192d7c768c.fdb0e26b708f7757c927665c32014118d81c9a82.zh-cn.xlf", "", "", "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.5b909c6e9d6b4335b3af5e7cc35d338fdda5ef7d.zh-cn.xlf")
$null = $ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/91dfa5fecd8eb92d70a2f5fb7a12f3cf15aa0586/e2e/e3b61664-96dc-4ab8-bc89-9c0d7fefc835.md", "", "", "0d2067ae-2ca4-448f-bc8e-89192d7c768c.md")
$null = $ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0f75b98accff473613980d6d00c49f8e386eb314/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e3b61664-96dc-4ab8-bc89-9c0d7fefc835.5b909c6e9d6b4335b3af5e7cc35d338fdda5ef7d.zh-cn.xlf", "", "", "0d2067ae-2ca4-448f-bc8e-89192d7c768c.fdb0e26b708f7757c927665c32014118d81c9a82.zh-cn.xlf")
$null = $ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/30fdea63b4e1cf4e0d65702035e6c8e5d3bb2d6f/e2e/e3b61664-96dc-4ab8-bc89-9c0d7fefc835.md", "", "", "0d2067ae-2ca4-448f-bc8e-89192d7c768c.md")
$null = $ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3436c8a9e40370fec3686ff289c357d4498874c7/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e3b61664-96dc-4ab8-bc89-9c0d7fefc835.5b909c6e9d6b4335b3af5e7cc35d338fdda5ef7d.zh-cn.xlf", "", "", "0d2067ae-2ca4-448f-bc8e-89192d7c768c.fdb0e26b708f7757c927665c32014118d81c9a82.zh-cn.xlf")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.5b909c6e9d6b4335b3af5e7cc35d338fdda5ef7d.de-de.xlf"
$ws.Range("E2").Value = "2016-03-23 18:55:52"
$ws.Range("F2").Value = "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.md"
$ws.Range("G2").Value = "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.5b909c6e9d6b4335b3af5e7cc35d338fdda5ef7d.de-de.xlf"
$ws.Range("H2").Value = "2016-03-23 18:56:45"
$ws.Range("J2").Value = "Include"

$ws.Range("A3").Value = "0d2067ae-2ca4-448f-bc8e-89192d7c768c.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "0d2067ae-2ca4-448f-bc8e-89192d7c768c.fdb0e26b708f7757c927665c32014118d81c9a82.de-de.xlf"
$ws.Range("E3").Value = "2016-03-23 18:57:33"
$ws.Range("F3").Value = "0d2067ae-2ca4-448f-bc8e-89192d7c768c.md"
$ws.Range("G3").Value = "0d2067ae-2ca4-448f-bc8e-89192d7c768c.fdb0e26b708f7757c927665c32014118d81c9a82.de-de.xlf"
$ws.Range("H3").Value = "2016-03-23 18:56:45"
$ws.Range("J3").Value = "Include"

$ws.Hyperlinks.Delete()
$null = $ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/91dfa5fecd8eb92d70a2f5fb7a12f3cf15aa0586/e2e/0d2067ae-2ca4-448f-bc8e-89192d7c768c.md", "", "", "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.md")
$null = $ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cfd74deb0f2683a3fdf3da7265b4c3859c833f76/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0d2067ae-2ca4-448f-bc8e-89192d7c768c.fdb0e26b708f7757c927665c32014118d81c9a82.de-de.xlf", "", "", "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.5b909c6e9d6b4335b3af5e7cc35d338fdda5ef7d.de-de.xlf")
$null = $ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/4de9de7fdcccb267520dd644a0c0017c94b76d54/e2e/0d2067ae-2ca4-448f-bc8e-89192d7c768c.md", "", "", "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.md")
$null = $ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/7dfb35f40aea915747a2a1d008ebcfacd58ad633/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0d2067ae-2ca4-448f-bc8e-89192d7c768c.fdb0e26b708f7757c927665c32014118d81c9a82.de-de.xlf", "", "", "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.5b909c6e9d6b4335b3af5e7cc35d338fdda5ef7d.de-de.xlf")
$null = $ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/91dfa5fecd8eb92d70a2f5fb7a12f3cf15aa0586/e2e/e3b61664-96dc-4ab8-bc89-9c0d7fefc835.md", "", "", "0d2067ae-2ca4-448f-bc8e-89192d7c768c.md")
$null = $ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cfd74deb0f2683a3fdf3da7265b4c3859c833f76/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e3b61664-96dc-4ab8-bc89-9c0d7fefc835.5b909c6e9d6b4335b3af5e7cc35d338fdda5ef7d.de-de.xlf", "", "", "0d2067ae-2ca4-448f-bc8e-89192d7c768c.fdb0e26b708f7757c927665c32014118d81c9a82.de-de.xlf")
$null = $ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/4de9de7fdcccb267520dd644a0c0017c94b76d54/e2e/e3b61664-96dc-4ab8-bc89-9c0d7fefc835.md", "", "", "0d2067ae-2ca4-448f-bc8e-89192d7c768c.md")
$null = $ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/7dfb35f40aea915747a2a1d008ebcfacd58ad633/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e3b61664-96dc-4ab8-bc89-9c0d7fefc835.5b909c6e9d6b4335b3af5e7cc35d338fdda5ef7d.de-de.xlf", "", "", "0d2067ae-2ca4-448f-bc8e-89192d7c768c.fdb0e26b708f7757c927665c32014118d81c9a82.de-de.xlf")

$wb.Save()
